$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.792813897132874
$ws.Range("B1").Value = 1.962782859802246
$ws.Range("C1").Value = 2.213757753372192
$ws.Range("D1").Value = 3.458035230636597
$ws.Range("E1").Value = 2.030245304107666
